# feat: add single and multi corrector
# Normalizes the lab-report sheet: strips stray leading index digits/markers
# that got prepended to the item-name column, fixes unit notation (10^9/L,
# 10^12/L, fL), and removes unit suffixes that had been baked into the
# result/reference-range values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "白细胞"
$ws.Range("C2").Value  = "10^9/L"

$ws.Range("A3").Value  = "红细胞"
$ws.Range("C3").Value  = "10^12/L"

$ws.Range("A4").Value  = "血红蛋白"

$ws.Range("A5").Value  = "血小板"
$ws.Range("C5").Value  = "10^9/L"

$ws.Range("A6").Value  = "红细胞压积"
# B6 used to hold "34.5%"; strip the percent sign but keep it text so
# Excel doesn't silently round-trip it through a float (e.g. 34.5 ->
# 34.500000000000004) when it's re-serialized.
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value  = "34.5"

$ws.Range("A7").Value  = "平均红细胞体积"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value  = "77.7"

$ws.Range("A8").Value  = "平均血红蛋白量"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value  = "26.4"

$ws.Range("A9").Value  = "平均血红蛋白浓度"

$ws.Range("A10").Value = "淋巴细胞百分率"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "17.4"

$ws.Range("A11").Value = "红细胞分布宽度"
$ws.Range("C11").Value = "fL"

$ws.Range("A12").Value = "红细胞分布宽度"

$ws.Range("A13").Value = "血小板分布宽度"

$ws.Range("A14").Value = ""

$ws.Range("A15").Value = "中性粒细胞百分率"

$ws.Range("A16").Value = "："

$ws.Range("A17").Value = "嗜酸性粒细胞百分率"
$ws.Range("D17").Value = "%"

$ws.Range("A18").Value = "嗜碱性粒细胞百分率"
$ws.Range("D18").Value = "%"

$ws.Range("A19").Value = "淋巴细胞绝对值"
$ws.Range("E19").Value = "1-3"

$ws.Range("A20").Value = "单"
$ws.Range("E20").Value = "0.07-0.33"

$ws.Range("A21").Value = "口"
$ws.Range("E21").Value = "3-5"

$ws.Range("A22").Value = ""
$ws.Range("E22").Value = "0.05-0.5"

$ws.Range("A23").Value = "口"
$ws.Range("E23").Value = "0.02-0.05"
